$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before row 1022. This shifts the existing
# rows 1022..1060 down to 1024..1062, and also pushes the current
# dimension/used range down accordingly (matching dimension A1:R1060 -> A1:R1062).
$ws.Range("A1022:A1023").EntireRow.Insert()

# --- New row 1022: Zafiro rojo, Primera, Provincia de Quillota ---
$ws.Cells.Item(1022, 1).Value2  = 4
$ws.Cells.Item(1022, 2).Value2  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(1022, 3).Value2  = "Los Lagos"
$ws.Cells.Item(1022, 4).Value2  = 45008
$ws.Cells.Item(1022, 5).Value2  = 10
$ws.Cells.Item(1022, 6).Value2  = 100112002
$ws.Cells.Item(1022, 7).Value2  = "Pimiento"
$ws.Cells.Item(1022, 8).Value2  = "Zafiro rojo"
$ws.Cells.Item(1022, 9).Value2  = "Primera"
$ws.Cells.Item(1022, 10).Value2 = 90
$ws.Cells.Item(1022, 11).Value2 = 29000
$ws.Cells.Item(1022, 12).Value2 = 29000
$ws.Cells.Item(1022, 13).Value2 = 29000
$ws.Cells.Item(1022, 14).Value2 = "$/caja 18 kilos"
$ws.Cells.Item(1022, 15).Value2 = "Provincia de Quillota"
$ws.Cells.Item(1022, 16).Value2 = 1611
$ws.Cells.Item(1022, 17).Value2 = 18
$ws.Cells.Item(1022, 18).Value2 = "Hortaliza"

# --- New row 1023: Zafiro verde, Primera, Provincia de Quillota ---
$ws.Cells.Item(1023, 1).Value2  = 4
$ws.Cells.Item(1023, 2).Value2  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(1023, 3).Value2  = "Los Lagos"
$ws.Cells.Item(1023, 4).Value2  = 45008
$ws.Cells.Item(1023, 5).Value2  = 10
$ws.Cells.Item(1023, 6).Value2  = 100112002
$ws.Cells.Item(1023, 7).Value2  = "Pimiento"
$ws.Cells.Item(1023, 8).Value2  = "Zafiro verde"
$ws.Cells.Item(1023, 9).Value2  = "Primera"
$ws.Cells.Item(1023, 10).Value2 = 90
$ws.Cells.Item(1023, 11).Value2 = 23000
$ws.Cells.Item(1023, 12).Value2 = 23000
$ws.Cells.Item(1023, 13).Value2 = 23000
$ws.Cells.Item(1023, 14).Value2 = "$/caja 18 kilos"
$ws.Cells.Item(1023, 15).Value2 = "Provincia de Quillota"
$ws.Cells.Item(1023, 16).Value2 = 1278
$ws.Cells.Item(1023, 17).Value2 = 18
$ws.Cells.Item(1023, 18).Value2 = "Hortaliza"

# Make sure the date cells keep the same date/time number format used
# throughout column D.
$ws.Range("D1022:D1023").NumberFormat = "YYYY-MM-DD HH:MM:SS"
